# "Added two new Mac-Addresses"
# Appends 10 new device rows (147-156) to the master-reg_center_device
# sheet, following the same pattern as the existing rows: regcntr_id
# 10001, sequential device_id starting at 3000166, lang_code "eng",
# is_active TRUE, cr_by "superadmin", cr_dtimes "now()".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

$startRow = 147
$startDeviceId = 3000166
$rowCount = 10

for ($i = 0; $i -lt $rowCount; $i++) {
    $row = $startRow + $i

    $ws.Cells.Item($row, 1).Value = 10001
    $ws.Cells.Item($row, 2).Value = ($startDeviceId + $i)
    $ws.Cells.Item($row, 3).Value = "eng"
    $ws.Cells.Item($row, 4).Value = $true
    $ws.Cells.Item($row, 5).Value = "superadmin"
    $ws.Cells.Item($row, 6).Value = "now()"
}

# Scroll/select to mirror where the author's cursor ended up after typing
# the new rows.
$excel.ActiveWindow.ScrollRow = 140
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("C152").Select()
